$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 475.33334
$ws.Range("I15").Value = 475.33334
$ws.Range("K15").Value = 1426.00002
$ws.Range("M15").Value = -1257.00002
# Row 17
$ws.Range("H17").Value = 2592.1428
$ws.Range("J17").Value = 2592.1428
$ws.Range("L17").Value = 7776.428400000001
$ws.Range("N17").Value = -8112.428400000001
# Row 33
$ws.Range("H33").Value = 177.5
$ws.Range("I33").Value = 177.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 177.5
$ws.Range("L33").Value = 0
$ws.Range("N33").Value = 51.5
$ws.Range("M33").ClearContents()
# Row 74
$ws.Range("H74").Value = 15189
$ws.Range("I74").Value = 3428
$ws.Range("K74").Value = 3428
$ws.Range("M74").Value = -2492
# Row 77
$ws.Range("H77").Value = 15189
$ws.Range("I77").Value = 3428
$ws.Range("K77").Value = 17140
$ws.Range("M77").Value = -12460
# Row 92
$ws.Range("H92").Value = 356.5
$ws.Range("I92").Value = 277.8
$ws.Range("J92").Value = 750
$ws.Range("K92").Value = 277.8
$ws.Range("L92").Value = 750
$ws.Range("M92").Value = 970.2
$ws.Range("N92").Value = -3246
# Row 100
$ws.Range("H100").Value = 2217.5
$ws.Range("I100").Value = 2461.8
$ws.Range("K100").Value = 2461.8
$ws.Range("M100").Value = -1920.8
# Row 135
$ws.Range("H135").Value = 880.5714
$ws.Range("I135").Value = 880.5714
$ws.Range("K135").Value = 7925.1426
$ws.Range("M135").Value = -5390.1426

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 963
$ws.Range("I32").Value = 907.8
$ws.Range("K32").Value = 907.8
$ws.Range("M32").Value = -620.8
# Row 61
$ws.Range("H61").Value = 2916.6667
$ws.Range("I61").Value = 2916.6667
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2916.6667
$ws.Range("L61").Value = 0
$ws.Range("N61").Value = -2704.6667
$ws.Range("M61").ClearContents()
# Row 136
$ws.Range("H136").Value = 2916.6667
$ws.Range("I136").Value = 2916.6667
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8750.000100000001
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = -6200.000100000001
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 548.5
$ws.Range("I80").Value = 284
$ws.Range("J80").Value = 1209.75
$ws.Range("K80").Value = 284
$ws.Range("L80").Value = 1209.75
$ws.Range("M80").Value = 714
$ws.Range("N80").Value = -3205.75
# Row 83
$ws.Range("H83").Value = 548.5
$ws.Range("I83").Value = 284
$ws.Range("J83").Value = 1209.75
$ws.Range("K83").Value = 1420
$ws.Range("L83").Value = 6048.75
$ws.Range("M83").Value = 3572
$ws.Range("N83").Value = -16032.75
# Row 86
$ws.Range("H86").Value = 6080.1113
$ws.Range("I86").Value = 1907.3334
$ws.Range("K86").Value = 1907.3334
$ws.Range("M86").Value = -784.3334
# Row 89
$ws.Range("H89").Value = 6080.1113
$ws.Range("I89").Value = 1907.3334
$ws.Range("K89").Value = 9536.666999999999
$ws.Range("M89").Value = -3920.666999999999
# Row 94
$ws.Range("H94").Value = 1050
$ws.Range("I94").Value = 99
$ws.Range("J94").Value = 2001
$ws.Range("K94").Value = 99
$ws.Range("L94").Value = 2001
$ws.Range("M94").Value = 352
$ws.Range("N94").Value = -2903
# Row 99
$ws.Range("H99").Value = 1588.8
$ws.Range("I99").Value = 986.25
$ws.Range("K99").Value = 986.25
$ws.Range("M99").Value = 511.75

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 299.5
$ws.Range("I2").Value = 4
$ws.Range("K2").Value = 4
$ws.Range("M2").Value = 109
# Row 4
$ws.Range("H4").Value = 29999992
$ws.Range("I4").Value = 29999990
$ws.Range("K4").Value = 29999990
$ws.Range("M4").Value = -29999878
# Row 68
$ws.Range("H68").Value = 46951.2
$ws.Range("J68").Value = 46951.2
$ws.Range("L68").Value = 46951.2
$ws.Range("N68").Value = -48449.2
# Row 71
$ws.Range("H71").Value = 46951.2
$ws.Range("J71").Value = 46951.2
$ws.Range("L71").Value = 140853.6
$ws.Range("N71").Value = -148341.6
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("N92").Value = 0
$ws.Range("L92").ClearContents()
# Row 99
$ws.Range("H99").Value = 2473.5
$ws.Range("J99").Value = 3000
$ws.Range("L99").Value = 3000
$ws.Range("N99").Value = -5996
# Row 126
$ws.Range("H126").Value = 2473.5
$ws.Range("J126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1228.2
$ws.Range("I5").Value = 984
$ws.Range("J5").Value = 1798
$ws.Range("K5").Value = 2952
$ws.Range("L5").Value = 5394
$ws.Range("M5").Value = -2840
$ws.Range("N5").Value = -5618
# Row 9
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("N70").Value = 0
$ws.Range("L70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("N73").Value = 0
$ws.Range("L73").ClearContents()
# Row 92
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("N92").Value = 0
$ws.Range("L92:M92").ClearContents()
# Row 109
$ws.Range("H109").Value = 958203.25
$ws.Range("I109").Value = 1429504.9
$ws.Range("K109").Value = 4288514.699999999
$ws.Range("M109").Value = -4287474.699999999
# Row 131
$ws.Range("H131").Value = 1578.762
$ws.Range("I131").Value = 607.8333
$ws.Range("K131").Value = 1823.4999
$ws.Range("M131").Value = 3216.5001
# Row 135
$ws.Range("H135").Value = 1228.2
$ws.Range("I135").Value = 984
$ws.Range("J135").Value = 1798
$ws.Range("K135").Value = 8856
$ws.Range("L135").Value = 16182
$ws.Range("M135").Value = -6321
$ws.Range("N135").Value = -21252

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 2000
$ws.Range("I70").Value = 2000
$ws.Range("K70").Value = 2000
$ws.Range("M70").Value = -1730
# Row 73
$ws.Range("H73").Value = 2000
$ws.Range("I73").Value = 2000
$ws.Range("K73").Value = 2000
$ws.Range("M73").Value = -1064
# Row 132
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1606.1111
$ws.Range("I22").Value = 1969.25
$ws.Range("J22").Value = 1315.6
$ws.Range("K22").Value = 1969.25
$ws.Range("L22").Value = 1315.6
$ws.Range("M22").Value = -1674.25
$ws.Range("N22").Value = -1905.6
# Row 27
$ws.Range("H27").Value = 1606.1111
$ws.Range("I27").Value = 1969.25
$ws.Range("J27").Value = 1315.6
$ws.Range("K27").Value = 1969.25
$ws.Range("L27").Value = 1315.6
$ws.Range("M27").Value = -1862.25
$ws.Range("N27").Value = -1529.6
# Row 40
$ws.Range("H40").Value = 5657.1055
$ws.Range("I40").Value = 5811
$ws.Range("K40").Value = 5811
$ws.Range("M40").Value = -5675
# Row 55
$ws.Range("H55").Value = 827.8
$ws.Range("J55").Value = 819.44446
$ws.Range("L55").Value = 819.44446
$ws.Range("N55").Value = -1165.44446
# Row 68
$ws.Range("H68").Value = 7142.857
$ws.Range("J68").Value = 7142.857
$ws.Range("L68").Value = 7142.857
$ws.Range("N68").Value = -8640.857
# Row 71
$ws.Range("H71").Value = 7142.857
$ws.Range("J71").Value = 7142.857
$ws.Range("L71").Value = 35714.285
$ws.Range("N71").Value = -43202.285
# Row 75
$ws.Range("H75").Value = 2107
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 78
$ws.Range("H78").Value = 2107
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# Row 100
$ws.Range("H100").Value = 8038.5557
$ws.Range("J100").Value = 9999.571
$ws.Range("L100").Value = 9999.571
$ws.Range("N100").Value = -11081.571
# Row 132
$ws.Range("H132").Value = 2050
$ws.Range("I132").Value = 2050
$ws.Range("K132").Value = 6150
$ws.Range("M132").Value = -3620

$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("L7").ClearContents()
# Row 9
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
# Row 49
$ws.Range("H49").Value = 15000
$ws.Range("I49").Value = 15000
$ws.Range("J49").Value = 15000
$ws.Range("K49").Value = 15000
$ws.Range("L49").Value = 15000
$ws.Range("M49").Value = -14770
$ws.Range("N49").Value = -15460
# Row 96
$ws.Range("H96").Value = 1072.25
$ws.Range("J96").Value = 944.5
$ws.Range("L96").Value = 944.5
$ws.Range("N96").Value = -3690.5
# Row 107
$ws.Range("H107").Value = 999
$ws.Range("J107").Value = 998
$ws.Range("L107").Value = 2994
$ws.Range("N107").Value = -6834
